$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    # MatchWholeWord=$true avoids accidental substring collisions
    # (e.g. "3.48" inside "-13.48").
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Row: Overall ATT (Group aggregation)
Replace-Text "3.18" "3.22"
Replace-Text "-10.99" "-11.07"
Replace-Text "1.47" "1.55"

# Row: Overall ATT (Dynamic aggregation)
Replace-Text "4.55" "4.42"
Replace-Text "-15.31" "-15.05"
Replace-Text "2.54" "2.27"

# Row: Cohort 2011
Replace-Text "6.31" "5.55"
Replace-Text "-21.24" "-19.76"
Replace-Text "3.48" "1.99"

# Row: Cohort 2014
Replace-Text "3.97" "3.56"
Replace-Text "-6.55" "-5.74"
Replace-Text "9.02" "8.21"

# Row: Cohort 2015
Replace-Text "3.50" "3.51"
Replace-Text "-13.48" "-13.49"
Replace-Text "0.23" "0.25"
